# Automatische test-sync: 2025-08-05 16:39:50
#
# Adds the new "Planning / Afspraak" test-mail row to the Logs sheet,
# bumps its matching Dashboard rollup row, extends the conditional
# formatting ranges and the chart's category/value source ranges to
# cover the new row.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dash = $wb.Worksheets.Item("Dashboard")

# --- Logs!A7:J7 - new mail log entry -----------------------------------
$logs.Range("A7").Value = "Kun jij dit even regelen?"
$logs.Range("B7").Value = "mailmind.test@zohomail.eu"
$logs.Range("C7").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D7").Value = "Planning / Afspraak"
$logs.Range("E7").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Range("F7").Value = "2025-08-05 16:39:02"
$logs.Range("G7").Value = "Ja"
$logs.Range("H7").Value = "Ja"
$logs.Range("I7").Value = "Nee"
$logs.Range("J7").Value = "Nee"

# --- Extend conditional formatting ranges on Logs to include row 7 ----
$dRules = $logs.Range("D2:D6").FormatConditions
for ($i = 1; $i -le $dRules.Count; $i++) {
    $dRules.Item($i).ModifyAppliesToRange($logs.Range("D2:D7"))
}

$gRules = $logs.Range("G2:G6").FormatConditions
for ($i = 1; $i -le $gRules.Count; $i++) {
    $gRules.Item($i).ModifyAppliesToRange($logs.Range("G2:G7"))
}

$hRules = $logs.Range("H2:H6").FormatConditions
for ($i = 1; $i -le $hRules.Count; $i++) {
    $hRules.Item($i).ModifyAppliesToRange($logs.Range("H2:H7"))
}

$iRules = $logs.Range("I2:I6").FormatConditions
for ($i = 1; $i -le $iRules.Count; $i++) {
    $iRules.Item($i).ModifyAppliesToRange($logs.Range("I2:I7"))
}

$jRules = $logs.Range("J2:J6").FormatConditions
for ($i = 1; $i -le $jRules.Count; $i++) {
    $jRules.Item($i).ModifyAppliesToRange($logs.Range("J2:J7"))
}

# --- Dashboard!A5:B5 - rollup row for the new category -----------------
$dash.Range("A5").Value = "Planning / Afspraak"
$dash.Range("B5").Value = 1

# --- Update chart source ranges to include the new Dashboard row ------
$chartObj = $dash.ChartObjects(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection(1)
$series.XValues = "='Dashboard'!`$A`$2:`$A`$5"
$series.Values = "='Dashboard'!`$B`$2:`$B`$5"
